$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'81.917.39"
$ws.Range('D2').Style = 'Normal'
$ws.Range("E2").Value = "'  +2.61%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range("D3").Value = "'3.163.82"
$ws.Range('D3').Style = 'Normal'
$ws.Range("E3").Value = "'  -1.29%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range("D5").Value = "'216.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range("E5").Value = "'  +4.94%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range("D6").Value = "'617.29"
$ws.Range('D6').Style = 'Normal'
$ws.Range("E6").Value = "'  -3.01%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range("E7").Value = "'  +17.54%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range("D8").Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range("D9").Value = "'0.581"
$ws.Range('D9').Style = 'Normal'
$ws.Range("E9").Value = "'  -4.90%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range("D10").Value = "'3.160.18"
$ws.Range('D10').Style = 'Normal'
$ws.Range("E10").Value = "'  -1.31%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range("D11").Value = "'0.601"
$ws.Range('D11').Style = 'Normal'
$ws.Range("E11").Value = "'  -3.17%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range("D12").Value = "'0.0000256"
$ws.Range('D12').Style = 'Normal'
$ws.Range("E12").Value = "'  +2.34%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range("B14").Value = "'Toncoin"
$ws.Range('B14').Style = 'Normal'
$ws.Range("C14").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C14').Style = 'Normal'
$ws.Range("D14").Value = "'5.29"
$ws.Range('D14').Style = 'Normal'
$ws.Range("E14").Value = "'  -2.89%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range("B15").Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B15').Style = 'Normal'
$ws.Range("C15").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C15').Style = 'Normal'
$ws.Range("D15").Value = "'3.743.81"
$ws.Range('D15').Style = 'Normal'
$ws.Range("E15").Value = "'  -1.20%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range("D16").Value = "'31.94"
$ws.Range('D16').Style = 'Normal'
$ws.Range("E16").Value = "'  -1.64%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range("D17").Value = "'81.597.67"
$ws.Range('D17').Style = 'Normal'
$ws.Range("E17").Value = "'  +2.43%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range("D18").Value = "'3.155.88"
$ws.Range('D18').Style = 'Normal'
$ws.Range("E18").Value = "'  -1.36%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range("D19").Value = "'3.21"
$ws.Range('D19').Style = 'Normal'
$ws.Range("E19").Value = "'  +8.35%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range("D20").Value = "'13.98"
$ws.Range('D20').Style = 'Normal'
$ws.Range("E20").Value = "'  -4.61%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range("D21").Value = "'434.07"
$ws.Range('D21').Style = 'Normal'
$ws.Range("E21").Value = "'  -0.23%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range("D22").Value = "'8.89"
$ws.Range('D22').Style = 'Normal'
$ws.Range("E22").Value = "'  -6.45%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range("D23").Value = "'5.11"
$ws.Range('D23').Style = 'Normal'
$ws.Range("E23").Value = "'  -2.49%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range("D24").Value = "'7.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range("E24").Value = "'  +5.67%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range("E25").Value = "'  +8.07%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range("D26").Value = "'11.81"
$ws.Range('D26').Style = 'Normal'
$ws.Range("E26").Value = "'  +7.15%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range("D27").Value = "'3.333.66"
$ws.Range('D27').Style = 'Normal'
$ws.Range("E27").Value = "'  -0.94%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range("D28").Value = "'76.46"
$ws.Range('D28').Style = 'Normal'
$ws.Range("E28").Value = "'  -1.75%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range("D29").Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range("E29").Value = "'  -0.48%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range("D30").Value = "'0.0000121"
$ws.Range('D30').Style = 'Normal'
$ws.Range("E30").Value = "'  +1.19%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range("D31").Value = "'0.999"
$ws.Range('D31').Style = 'Normal'
$ws.Range("E31").Value = "'  -0.03%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range("E32").Value = "'  -3.40%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range("D33").Value = "'567.36"
$ws.Range('D33').Style = 'Normal'
$ws.Range("E33").Value = "'  +6.94%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range("E34").Value = "'  +0.44%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range("D35").Value = "'0.146"
$ws.Range('D35').Style = 'Normal'
$ws.Range("E35").Value = "'  +18.68%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range("E36").Value = "'  +5.20%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range("D37").Value = "'1.98"
$ws.Range('D37').Style = 'Normal'
$ws.Range("E37").Value = "'  -1.54%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range("D38").Value = "'22.60"
$ws.Range('D38').Style = 'Normal'
$ws.Range("E38").Value = "'  -3.11%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range("D39").Value = "'0.999"
$ws.Range('D39').Style = 'Normal'
$ws.Range("E39").Value = "'  +0.04%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range("D40").Value = "'6.08"
$ws.Range('D40').Style = 'Normal'
$ws.Range("E40").Value = "'  +9.40%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range("D41").Value = "'0.406"
$ws.Range('D41').Style = 'Normal'
$ws.Range("E41").Value = "'  -1.36%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range("E42").Value = "'  +3.96%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range("D43").Value = "'3.02"
$ws.Range('D43').Style = 'Normal'
$ws.Range("E43").Value = "'  +15.70%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range("D44").Value = "'2.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range("E44").Value = "'  +9.54%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range("D45").Value = "'158.68"
$ws.Range('D45').Style = 'Normal'
$ws.Range("E45").Value = "'  -3.74%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range("E46").Value = "'  +0.04%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range("D47").Value = "'186.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range("E47").Value = "'  -3.11%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range("D48").Value = "'44.73"
$ws.Range('D48').Style = 'Normal'
$ws.Range("E48").Value = "'  +2.88%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range("E49").Value = "'  +0.06%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range("D50").Value = "'26.47"
$ws.Range('D50').Style = 'Normal'
$ws.Range("E50").Value = "'  +2.51%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range("D51").Value = "'0.763"
$ws.Range('D51').Style = 'Normal'
$ws.Range("E51").Value = "'  -5.60%  "
$ws.Range('E51').Style = 'Normal'
